$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers: Latitude -> Lat, Longitude -> Lon
# (Lon is entered before Lat so the shared-string table order matches
# the authored workbook: "Lon" then "Lat".)
$ws.Range("C1").Value = "Lon"
$ws.Range("B1").Value = "Lat"

# Fill in lat/lon decimal-degree coordinates for all site rows (rows 6-12 previously blank)
$ws.Range("B6").Value = 27.9338
$ws.Range("C6").Value = 30.8836

$ws.Range("B7").Value = 30.5727
$ws.Range("C7").Value = 31.51

$ws.Range("B8").Value = 27.899
$ws.Range("C8").Value = 30.8666

$ws.Range("B9").Value = 37.9891
$ws.Range("C9").Value = 23.7322

$ws.Range("B10").Value = 40.7514
$ws.Range("C10").Value = 14.485

$ws.Range("B11").Value = 55.1666
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 45.55
$ws.Range("C12").Value = 61.8666

# Update site name for row 11 (Northumbria, Ireland -> Northumberland, England)
$ws.Range("A11").Value = "Northumberland, England"

# Update the active selection to match the final saved state
[void]$ws.Range("D17").Select()
